# Applies the edit described by the commit "Fruta / hortaliza, semanal":
# a new daily price observation row is inserted at row 27 of the sheet,
# pushing the existing rows 27-146 down to 28-147.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 27 (shifts rows 27..146 -> 28..147,
# carrying their formatting/styles with them automatically).
$ws.Rows.Item(27).Insert()

# Populate the newly inserted row 27 with the new observation.
$ws.Cells.Item(27, 1).Value2  = 4
$ws.Cells.Item(27, 2).Value2  = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(27, 3).Value2  = "Los Lagos"
$ws.Cells.Item(27, 4).Value2  = 45076
$ws.Cells.Item(27, 5).Value2  = 10
$ws.Cells.Item(27, 6).Value2  = 100112031
$ws.Cells.Item(27, 7).Value2  = "Poroto verde"
$ws.Cells.Item(27, 8).Value2  = "Magnum"
$ws.Cells.Item(27, 9).Value2  = "Primera"
$ws.Cells.Item(27, 10).Value2 = 45
$ws.Cells.Item(27, 11).Value2 = 30000
$ws.Cells.Item(27, 12).Value2 = 30000
$ws.Cells.Item(27, 13).Value2 = 30000
$ws.Cells.Item(27, 14).Value2 = "`$/malla 25 kilos"
$ws.Cells.Item(27, 15).Value2 = "Perú"
$ws.Cells.Item(27, 16).Value2 = 1200
$ws.Cells.Item(27, 17).Value2 = 25
$ws.Cells.Item(27, 18).Value2 = "Hortaliza"
